$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'294.20"
$ws.Range("E2").Value = "'1.23%"
$ws.Range("D3").Value = "'31.03"
$ws.Range("E3").Value = "'0.70%"
$ws.Range("D4").Value = "'4.925"
$ws.Range("E4").Value = "'1.01%"
$ws.Range("D5").Value = "'0.07340"
$ws.Range("E5").Value = "'2.41%"
$ws.Range("D6").Value = "'2.295"
$ws.Range("E6").Value = "'29.63%"
$ws.Range("D7").Value = "'7.744"
$ws.Range("E7").Value = "'0.81%"
$ws.Range("D8").Value = "'3.751"
$ws.Range("E8").Value = "'0.36%"
$ws.Range("D9").Value = "'0.9083"
$ws.Range("E9").Value = "'1.57%"
$ws.Range("D10").Value = "'0.1688"
$ws.Range("E10").Value = "'1.31%"
$ws.Range("D11").Value = "'0.07979"
$ws.Range("E11").Value = "'7.22%"
$ws.Range("D12").Value = "'0.08147"
$ws.Range("E12").Value = "'0.66%"
$ws.Range("D13").Value = "'0.03100"
$ws.Range("E13").Value = "'4.20%"
$ws.Range("E14").Value = "'0.71%"
$ws.Range("D15").Value = "'0.001511"
$ws.Range("E15").Value = "'1.13%"
$ws.Range("D16").Value = "'0.005852"
$ws.Range("E16").Value = "'-0.20%"
$ws.Range("D17").Value = "'3.485"
$ws.Range("E17").Value = "'0.80%"
$ws.Range("D18").Value = "'2.076"
$ws.Range("E18").Value = "'-1.34%"
$ws.Range("D19").Value = "'0.3328"
$ws.Range("E19").Value = "'1.05%"
$ws.Range("E20").Value = "'0.42%"
$ws.Range("D21").Value = "'3.968"
$ws.Range("E21").Value = "'-9.46%"
$ws.Range("E22").Value = "'4.92%"
$ws.Range("D23").Value = "'0.04544"
$ws.Range("E23").Value = "'1.65%"
$ws.Range("E24").Value = "'-0.31%"
$ws.Range("D25").Value = "'0.004647"
$ws.Range("E25").Value = "'15.75%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'3.96%"
$ws.Range("D27").Value = "'0.0003397"
$ws.Range("D39").Value = "'0.01603"
$ws.Range("E39").Value = "'-2.86%"
$ws.Range("D40").Value = "'0.04441"
$ws.Range("E40").Value = "'2.27%"
$ws.Range("D41").Value = "'0.007340"
$ws.Range("E41").Value = "'-0.98%"
$ws.Range("D42").Value = "'0.1328"
$ws.Range("E42").Value = "'1.38%"
$ws.Range("D43").Value = "'0.008631"
$ws.Range("D44").Value = "'0.002023"
$ws.Range("E44").Value = "'-0.82%"
$ws.Range("D45").Value = "'0.009518"
$ws.Range("E45").Value = "'-6.41%"
$ws.Range("D46").Value = "'0.00005970"
$ws.Range("E46").Value = "'4.21%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("E48").Value = "'2.38%"
$ws.Range("E49").Value = "'-3.48%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("E51").Value = "'-0.05%"
